$d = $word.ActiveDocument

function Replace-RunText($findText, $newText) {
    # Find the target run's text in the document.
    $rng = $d.Content
    $rng.Find.Text = $findText
    $rng.Find.Execute() | Out-Null

    # Pull the run's opening markup (including w:rPr) straight out of the
    # live OOXML so the replacement run keeps identical run formatting.
    # The text inside WordOpenXML is itself XML-escaped, so search using
    # an escaped copy of the needle.
    $xml = $rng.WordOpenXML
    $findEscaped = $findText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $idx = $xml.IndexOf($findEscaped)
    $head = $xml.Substring(0, $idx)
    $rStart = $head.LastIndexOf("<w:r>")
    if ($rStart -lt 0) {
        $rStart = $head.LastIndexOf("<w:r ")
    }
    $tTagStart = $head.LastIndexOf("<w:t")
    $runPrefix = $xml.Substring($rStart, $tTagStart - $rStart)

    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $runPrefix + '<w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($frag)
    Write-Output "inserted replacement run"

    # The original text is still present right before what we just
    # inserted; find it again and remove it.
    $old = $d.Content
    $old.Find.Text = $findText
    $old.Find.Execute() | Out-Null
    $old.Delete()
    Write-Output "removed original run"
}

Replace-RunText "Verantwortlich für den Aufbau von SmartCash und unterstützenden Anwendungen." "Verantwortlich für die Softwareentwicklung von SmartCash."

Replace-RunText "Dieses Team konzentriert sich auf Gemeinschaftsbildung, Wachstum und allgemeine Nutzerakquise in Südamerika" "Dieses Team richtet seinen Fokus auf Communitybildung, Wachstum und dem Anwerben neuer Mitglieder in Südamerika"

Replace-RunText "Dieser Hive ist verantwortlich für On-Boarding & allgemeinen SmartCash Support." "Dieses Team ist für das On-Boarding und generelle SmartCash Unterstützung verantwortlich."

Replace-RunText "MITARBEIT GEFÄLLIG?" "WILLST DU MITMACHEN?"

Replace-RunText "SmartHive hat Platz für Leute aus allen Erfahrungsbereichen. Wissensdurst gefragt!" "SmartHive hat Platz für Leute aus allen Erfahrungsbereichen. Sei motiviert!"

Replace-RunText "Wir glauben, dass „Core“-Teams eine schlechte Idee sind, weil diese letztendlich zu Ineffizienz und Korruption führen. Wir wollen es überwinden und ein dezentralisiertes Organisationsmodell schaffen, das von Ameinse- und Bienenkolonien inspiriert ist." "Wir glauben, dass „Core“ Teams eine schlechte Idee sind und letztendlich zu Ineffizienz und Korruption führen. Wir wollen dies hinter uns lassen und ein dezentralisiertes Organisationsmodell schaffen, das von Ameisen- und Bienenvölkern inspiriert ist."

Replace-RunText "Um eine dezentrale Leitungsstruktur zu realisieren und zu erhalten, führen wir zwei Konzepte ein: SmartHive und Hive Strukturierungs Teams (HST). SmartHive ermöglicht jedem der Coins hält die Gelegenheit über Vorschläge abzustimmen, die von der Gemeinschaft eingereicht werden. SmartHive wird das Lebenselixier des Projekts sein, welches jedem ermöglichen soll, sich einzubringen und Vorschläge einreichen zu können. Mit dieser Hilfe kann organisches Wachstum von der Basis aus erzeugt werden, was eine Management Struktur zur Folge hat, die von unten nach oben gerichtet ist." "Um eine dezentrale Entscheidungsstruktur zu schaffen und zu erhalten, führen wir zwei Konzepte ein: SmartHive und Hive Strukturierungs Teams (HST). SmartHive ermöglicht es jedem, der SmartCash besitzt, über von der Community eingereichte Anträge, abzustimmen. SmartHive wird das Lebenselixier des Projekts sein, welches jedem ermöglichen soll, sich einzubringen und Vorschläge einreichen zu können. Mit dieser Hilfe kann organisches Wachstum von der Basis aus erzeugt werden, was eine Management Struktur zur Folge hat, die von unten nach oben gerichtet ist."
